# Applies the edits described by the diff to the active document.
# Uses $d.Content.Find.Execute for straightforward text replacements and
# $d.Paragraphs / Range based manipulation for structural changes
# (paragraph insertion, run restructuring for proofErr-wrapped words).

$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $range = $d.Content
    $ok = $range.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        Write-Output "MISSING: $findText"
    }
}

# 1. "And then the project followed..." -> "project followed..."
Replace-Text "And then the project followed a generic implementation and testing " "project followed a generic implementation and testing "

# 2. "We began by designing and implementing our own WSN protocol based on what we've already seen"
#    -> "started by designing and implementing our own WSN protocol "
Replace-Text "We began by designing and implementing our own WSN protocol based on what we’ve already seen" "started by designing and implementing our own WSN protocol "

# 3. "the network's parameters such that it had very little idle time" -> "the network's parameters "
Replace-Text " the network’s parameters such that it had very little idle time" " the network’s parameters "

# 4. "We produced response time models and evaluated them against the actual values"
#    -> "We produced response time models and evaluated them "
Replace-Text "We produced response time models and evaluated them against the actual values" "We produced response time models and evaluated them "

# 5. " how well the protocol handled external interference" -> " intergirty against external interference"
Replace-Text " how well the protocol handled external interference" " intergirty against external interference"

# 6. "In brief, a WSN is a large interconnected network" -> "A WSN is a large interconnected network"
Replace-Text "In brief, a WSN" "A WSN"

# 7. Full rewrite of the "Protocol needs to automatically handle..." paragraph
Replace-Text "Protocol needs to automatically handle the effects of external interference (coming from outside the network) – whether that’s intentionally created to affect the network or someone putting their lunch in the microwave and it accidentally causing an effect. and also prevent internal interference i.e. make sure the transmissions don’t overlap with eachother " "WSN protocol used must be able to manage external and internal interference. External – coming from outside the network (deliberate – trying to disrupt network, accidental – heating your lunch in a microwave nearby). Internal – within the network, transmissions overlapping with eachother. Both types need to be managed to ensure the network makes progress."

# 8. Full rewrite of the "There are a wide range of protocols..." paragraph
Replace-Text "There are a wide range of protocols to choose from, but the common ground is the use of the IEEE 802.15.4 physical layer because of it slow rate low power capabilities. It provides two modes, either beacon enabled or non-beacon enabled with a large range frequencies to transmit on – known as channels. " "There are a wide range of protocols to choose from, but the common ground is the use of the IEEE 802.15.4 physical layer because of it low power operation. It provides two operating modes, either beacon enabled or non-beacon enabled. It provides a large range of frequencies to transmit on, broken down into sections known as channels"

# 9. Full rewrite of the WHART/superframe paragraph
Replace-Text ") is essentially “transmit when you and the channel are ready”. BE – WirelessHART follows a completely deterministic synchornised ordering of transmission based on the assignment of timeslots (bandwidth) to transmissions which make up a superframe. This scheduling gets rid of the issue of overlapping transmissions as the devices know when they are allowed to transmit. External int is managed by channel hopping – changing channel each transmission and also blacklisitng certain frequencies when they are too noisy. " ") is essentially “transmit when you and the channel are ready”. BE e.g. WHART. Assigns each transmission a timeslot (a bit of bandwidth) which are then ordered to form a superframe which prevents transmissions from overlapping naturally – as long as all of the devices are synchronized (by the network manager node). External interference is managed by channel hopping and blacklisting – changing channels every transmission and avoiding channels that are too noisy"

# 11. "TinyOS is an OS which uses a C-like language" -> "TinyOS uses a C-like language"
Replace-Text " is an OS which uses a C-like language " " uses a C-like language "

# 12. "supports most strictly typed HLLs" -> "supports most HLLs"
Replace-Text ". Mote Runner which supports most strictly typed HLLs and " ". Mote Runner which supports most HLLs and "

# 14. "...we'd later prove. We introduced" -> "...we'd later prove by measurement. We introduced"
Replace-Text "we’d later prove. We introduced" "we’d later prove by measurement. We introduced"

# 15. Full rewrite of "For all of our analysis..." sentence start
Replace-Text "For all of our analysis we employed the use of the simulator and use the logger system library ( a debugging tool) to act as packet " "Simulator used for analysis and use the logger outputs were used to act as packet "

# 16. Remove "However these results need to be taken..." trailing sentence
Replace-Text " – when a device received a packet it outputted “packet received”. However these results need to be taken with a pinch of salt as the simulator is more of a software model of the motes rather than a hardware emulator – the actual physical motes may run differently. " " – when a device received a packet it outputted “packet received”. "

# 18. Rewrite the "We first began by optimizing the timeslot length (..." paragraph tail
Replace-Text "We first began by optimizing the timeslot length (" "Optimising the timeslot length ("
Replace-Text " assigned) to devices to complete a transmission. Larger payload = more time so we optimized this for a range of payload sizes based on the spread of the data transfers specification – we wanted a 0 packet fault rate (number of packets lost) for the shortest timeslot possible." " assigned) to devices to complete a transmission. Larger payload = more time so we optimized this for a range of payload sizes based on the spread of the data transfers specification. we wanted packet fault rate (number of packets lost) of 0 for the smallest timeslot possible."

# 19. Replace the "TO find the synchronisation constant..." paragraph with the new shorter sentence
Replace-Text "TO find the synchronisation constant we just used the logger to measure the time between the end of a superframe and a device recivieng the synchronisation broadcast and we’d do this over all the optimized timlsots to ensure it is in fact constant. " "Synchronisation constant found by measuring time between end of superframe and start of next."

# 20. Replace the "We tested the intergirty of LikeWHART..." paragraph
Replace-Text "We tested the intergirty of LikeWHART against external interference by introducing a rogue node that simply changed channel at a rate we’d increase in an attempt to see the effect on the packet fault rate" "Intergrity of LWHART tested by introducing a rogue device which cycled through each of the 16 channels broadcasting on them, We increased rate of channel changing to see the effects on network."

# 21. "Finally we computed the theoretical response times..." -> "theoretical response times calculated..."
Replace-Text "Finally we computed the theoretical response times for our transfer specification and then " "theoretical response times calculated for our transfer specification and then "

# 22. "We found that the optimal timeslot doesn't actually double..." -> "optimal timeslot doesn't actually double..."
Replace-Text "We found that the optimal timeslot doesn’t actually double" "optimal timeslot doesn’t actually double"

# 23. Rewrite trailing sentence about response times / measuring apparatus
Replace-Text " computed response times for the data transfer specification and they were the same with only limited fluctuations but we put this down to the measuring apparatus. " " computed response times for the data transfer specification – our response time models were right!"

# 24. "was hardly unsurprisingly" -> "was hardly surprisingly"
Replace-Text "Our integrity analysis was hardly unsurprisingly." "Our integrity analysis was hardly surprisingly."

# 25. append new sentence after "...is terrible at handling external interference"
Replace-Text " is terrible at handling external interference" " is terrible at handling external interference. It doesn’t manage it, it just tries to avoid it"

# 26. Fix "integrirty" -> "integrity" typo
Replace-Text "integrirty" "integrity"

Write-Output "done"
